$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to retain text formatting so numeric-looking
# strings (e.g. "5.92", "600.47") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.842.68'
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = '3.473.62'
$ws.Range("E3").Value = '  -2.47%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '600.47'
$ws.Range("E5").Value = '  -3.03%  '
$ws.Range("D6").Value = '147.52'
$ws.Range("E6").Value = '  -4.76%  '
$ws.Range("D7").Value = '3.472.04'
$ws.Range("E7").Value = '  -2.48%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("E11").Value = '  +3.88%  '
$ws.Range("E12").Value = '  -3.63%  '
$ws.Range("E13").Value = '  -4.38%  '
$ws.Range("D14").Value = '4.059.56'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("E15").Value = '  -6.70%  '
$ws.Range("D16").Value = '3.471.67'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '66.831.49'
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  -5.35%  '
$ws.Range("D20").Value = '10.17'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = '15.16'
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("D22").Value = '432.44'
$ws.Range("E22").Value = '  -4.80%  '
$ws.Range("E23").Value = '  -6.01%  '
$ws.Range("D24").Value = '79.16'
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = '3.609.99'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E27").Value = '  -9.49%  '
$ws.Range("D28").Value = '9.80'
$ws.Range("E28").Value = '  -6.76%  '
$ws.Range("E29").Value = '  -10.26%  '
$ws.Range("E30").Value = '  -3.32%  '
$ws.Range("D31").Value = '1.60'
$ws.Range("E31").Value = '  -6.59%  '
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").Value = '25.31'
$ws.Range("E34").Value = '  -3.08%  '
$ws.Range("D35").Value = '3.462.93'
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '5.92'
$ws.Range("E36").Value = '  -7.43%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  -6.25%  '
$ws.Range("D39").Value = '7.90'
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '172.37'
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("E42").Value = '  -3.65%  '
$ws.Range("D43").Value = '2.08'
$ws.Range("E43").Value = '  -12.91%  '
$ws.Range("E44").Value = '  -4.52%  '
$ws.Range("D45").Value = '0.897'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '46.38'
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '28.88'
$ws.Range("E47").Value = '  -7.42%  '
$ws.Range("E48").Value = '  -8.25%  '
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("E50").Value = '  -10.35%  '
$ws.Range("D51").Value = '0.968'
$ws.Range("E51").Value = '  -4.82%  '

# Restore default (Normal) style on column D so no stray number format
# remains attached to the cells (matches original workbook styling).
$ws.Range("D2:D51").Style = "Normal"

